# Update odds values on the active sheet (FlashScore weekly games export)
# to match the latest scrape, per the commit "Atualizando o arquivo XLSX".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.2
$ws.Range("I2").Value = 4.2
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 5.5
$ws.Range("M2").Value = 1.18
$ws.Range("N2").Value = 4.5
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 1.25
$ws.Range("S2").Value = 1.85
$ws.Range("T2").Value = 1.95
$ws.Range("U2").Value = 3
$ws.Range("V2").Value = 1.36
$ws.Range("X2").Value = 8
$ws.Range("AD2").Value = 6.5
$ws.Range("AG2").Value = 7
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 19
$ws.Range("AL2").Value = 81
$ws.Range("AO2").Value = 15
$ws.Range("AT2").Value = 1.83
$ws.Range("AU2").Value = 13
$ws.Range("AX2").Value = 34
$ws.Range("BA2").Value = 251

# Row 3
$ws.Range("Q3").Value = 2.05
$ws.Range("R3").Value = 1.75
$ws.Range("AC3").Value = 9.5

# Row 4
$ws.Range("N4").Value = 8.5
$ws.Range("Z4").Value = 11
$ws.Range("AC4").Value = 8.5
$ws.Range("AK4").Value = 51
$ws.Range("AX4").Value = 34

# Row 6
$ws.Range("G6").Value = 2.4
$ws.Range("AR6").Value = 51
$ws.Range("AW6").Value = 5

# Row 9
$ws.Range("G9").Value = 2.05
$ws.Range("I9").Value = 3.7
$ws.Range("J9").Value = 2.88
$ws.Range("N9").Value = 7.5
$ws.Range("U9").Value = 1.91
$ws.Range("V9").Value = 1.8
$ws.Range("Y9").Value = 9.5
$ws.Range("Z9").Value = 19
$ws.Range("AH9").Value = 17
$ws.Range("AI9").Value = 13
$ws.Range("AZ9").Value = 67
$ws.Range("BB9").Value = 251

# Row 10
$ws.Range("G10").Value = 1.42
$ws.Range("I10").Value = 7
$ws.Range("U10").Value = 2.2
$ws.Range("V10").Value = 1.62
$ws.Range("Z10").Value = 9
$ws.Range("AH10").Value = 34
$ws.Range("AI10").Value = 21
$ws.Range("AN10").Value = 3.25
$ws.Range("AU10").Value = 9.5
$ws.Range("AZ10").Value = 151

# Row 11
$ws.Range("O11").Value = 1.25
$ws.Range("P11").Value = 3.75
$ws.Range("Q11").Value = 1.83
$ws.Range("R11").Value = 2.03
